# httk-benchmarks.xlsx : "Added invivoPKfit outputs to dashboard script"
#
# A new benchmark row (version 2.3.0) is appended to Table1 on Sheet1,
# along with its two new notes/version shared strings, and the view is
# left scrolled/selected on the new Notes cell (R24), matching what a
# user does after typing a new row into the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by one row (this also keeps the table's `ref` /
# `autoFilter` range and the sheet `dimension` in sync automatically).
$newListRow = $lo.ListRows.Add()
$newRow = $newListRow.Index + 1   # header is row 1, so data row index -> sheet row

# Match the formatting (left-aligned, same as every other data row) of
# the existing rows before filling in the values.
$ws.Range("A" + $newRow + ":R" + $newRow).HorizontalAlignment = -4131

# New benchmark results row for release 2.3.0.
$ws.Range("A$newRow").Value = "2.3.0"
$ws.Range("B$newRow").Value = 1023
$ws.Range("C$newRow").Value = 0.9999
$ws.Range("D$newRow").Value = 1
$ws.Range("E$newRow").Value = 1
$ws.Range("F$newRow").Value = 1.063
$ws.Range("G$newRow").Value = 352
$ws.Range("H$newRow").Value = 0.2996
$ws.Range("I$newRow").Value = 352
$ws.Range("J$newRow").Value = 1.419
$ws.Range("K$newRow").Value = 86
$ws.Range("L$newRow").Value = 1.047
$ws.Range("M$newRow").Value = 86
$ws.Range("N$newRow").Value = 1.33
$ws.Range("O$newRow").Value = 86
$ws.Range("P$newRow").Value = 0.6344
$ws.Range("Q$newRow").Value = 863
$ws.Range("R$newRow").Value = "Used Caco-2 to replace Fabs=Fgut=1"

# Leave the sheet scrolled to / selecting the new Notes cell, as if the
# user had just finished typing the new row.
$ws.Activate()
$ws.Range("R" + $newRow).Select()

Write-Host "Added benchmark row $newRow (version 2.3.0) to Table1"
